# Update "想去人数" (want-to-go count) values in column F
# for the "展览" (sheet1) and "全部类型" (sheet4) worksheets.
# The two sheets list mostly the same events but "全部类型" has one
# extra row inserted near the top, so the target rows differ by sheet.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> ordered list of (Row, NewValue) pairs
$updates = @{
    "展览"   = @(
        @{Row = 3;  Value = 1853},
        @{Row = 4;  Value = 484},
        @{Row = 7;  Value = 2513},
        @{Row = 8;  Value = 160},
        @{Row = 9;  Value = 86},
        @{Row = 11; Value = 1514},
        @{Row = 12; Value = 523},
        @{Row = 15; Value = 226},
        @{Row = 19; Value = 218},
        @{Row = 20; Value = 215},
        @{Row = 22; Value = 158},
        @{Row = 24; Value = 1598},
        @{Row = 30; Value = 408}
    )
    "全部类型" = @(
        @{Row = 3;  Value = 1853},
        @{Row = 5;  Value = 484},
        @{Row = 8;  Value = 2513},
        @{Row = 9;  Value = 160},
        @{Row = 10; Value = 86},
        @{Row = 12; Value = 1514},
        @{Row = 13; Value = 523},
        @{Row = 16; Value = 226},
        @{Row = 20; Value = 218},
        @{Row = 21; Value = 215},
        @{Row = 23; Value = 158},
        @{Row = 25; Value = 1598},
        @{Row = 31; Value = 408}
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
